$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Ticket Code" (C) and "Scanned" (D) columns are being wiped out.
# Column C keeps an (empty, quote-prefixed) text value in every row, including
# the header. Column D's header also becomes an empty quote-prefixed text
# value, but D's data rows (2-10) become fully blank cells.

# Header row: C1/D1 -> empty quote-prefixed text (was "Ticket Code"/"Scanned").
$ws.Range("C1").Value = "'"
$ws.Range("D1").Value = "'"

# Data rows: C2:C10 -> empty quote-prefixed text (was per-row hash codes).
$ws.Range("C2:C10").Value = "'"

# Data rows: D2:D10 -> fully blank. D2 previously held an empty quote-prefixed
# text value, so first overwrite it with a plain value to drop the
# quote-prefix formatting, matching D3:D10, then clear all of D2:D10.
$ws.Range("D2").Value = "x"
$ws.Range("D2:D10").ClearContents()
